$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '34.809.36'
$ws.Cells.Item(2, 5).Value = '  -0.09%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.824.72'
$ws.Cells.Item(3, 5).Value = '  +0.91%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.30%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''230.39'
$ws.Cells.Item(5, 5).Value = '  -0.13%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''0.615'
$ws.Cells.Item(6, 5).Value = '  +1.14%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.29%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''39.71'
$ws.Cells.Item(8, 5).Value = '  +0.21%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.319'
$ws.Cells.Item(9, 5).Value = '  +3.19%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.0682'
$ws.Cells.Item(10, 5).Value = '  +0.15%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.86%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '2.091.98'
$ws.Cells.Item(12, 5).Value = '  +1.01%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''11.26'
$ws.Cells.Item(13, 5).Value = '  +2.20%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.827.36'
$ws.Cells.Item(14, 5).Value = '  +0.94%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.665'
$ws.Cells.Item(15, 5).Value = '  +1.77%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''4.63'
$ws.Cells.Item(16, 5).Value = '  -0.21%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '34.877.17'
$ws.Cells.Item(17, 5).Value = '  +0.11%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''69.38'
$ws.Cells.Item(18, 5).Value = '  +1.02%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '0.0₃0785'
$ws.Cells.Item(19, 5).Value = '  +0.49%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''239.30'
$ws.Cells.Item(20, 5).Value = '  +1.16%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''12.09'
$ws.Cells.Item(21, 5).Value = '  +3.28%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.17%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.08%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''2.25'
$ws.Cells.Item(24, 5).Value = '  -0.29%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''173.49'
$ws.Cells.Item(25, 5).Value = '  +0.33%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''7.70'
$ws.Cells.Item(26, 5).Value = '  -0.88%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''0.123'
$ws.Cells.Item(27, 5).Value = '  +3.04%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''17.27'
$ws.Cells.Item(28, 5).Value = '  -0.36%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -5.08%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.17%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''0.0548'
$ws.Cells.Item(31, 5).Value = '  +0.24%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''3.90'
$ws.Cells.Item(32, 5).Value = '  +1.01%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''3.91'
$ws.Cells.Item(33, 5).Value = '  -0.99%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(34, 4).Value = '''1.21'
$ws.Cells.Item(34, 5).Value = '  +3.91%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(35, 4).Value = '''1.81'
$ws.Cells.Item(35, 5).Value = '  +1.63%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +11.79%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''0.696'
$ws.Cells.Item(37, 5).Value = '  +3.31%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''91.78'
$ws.Cells.Item(38, 5).Value = '  -1.43%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '1.337.07'
$ws.Cells.Item(39, 5).Value = '  +2.83%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''1.01'
$ws.Cells.Item(40, 5).Value = '  +2.97%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +1.07%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''14.42'
$ws.Cells.Item(42, 5).Value = '  -1.58%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -0.69%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''2.24'
$ws.Cells.Item(44, 5).Value = '  -3.02%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.07%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.40%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +2.04%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '2.008.98'
$ws.Cells.Item(48, 5).Value = '  +1.19%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.21%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''0.0669'
$ws.Cells.Item(50, 5).Value = '  +4.14%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Quant'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(51, 4).Value = '''98.14'
$ws.Cells.Item(51, 5).Value = '  -0.49%  '
